# Disable "as you type" autoformatting so straight quotes in the source
# text are not silently converted to curly quotes while we edit.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    $r.Text = $newText
}

# 1) Title heading
Replace-ExactText "Review 168: [Short] VERA: VECTOR-BASED RANDOM MATRIX ADAPTATION" "Review 167: [Short] Reward-Augmented Decoding: Efficient Controlled Text Generation With a Unidirectional Reward Model"

# 2) Paper link (bold run)
Replace-ExactText "Paper: https://arxiv.org/abs/2310.11454v2" "Paper: https://arxiv.org/abs/2310.09520v4"

# 3) HuggingFace link
Replace-ExactText "https://huggingface.co/papers/2310.11454" "https://huggingface.co/papers/2310.09520"

# 4) Remove the paragraph containing the LoRA intro text, along with the two
#    preceding empty paragraphs (they become superfluous once the Hebrew
#    introduction is replaced by the new, longer article text further down).
$targetText = "כבר סקרנו השבוע מאמר שמציע שיפור ל-LoRA שיטת טיוב(finetune) מודלי שפה חסכונית מבחינת משאבי חישוב הנדרשים. היום ב-#shorthebrewpapereviews נסקור שיפור נוסף ל-LoRA המאפשר להקטין את כמות המשאבים הנדרשים לטיוב עוד יותר. "
$foundStart = -1
$foundEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {
        $foundStart = $i - 2
        $foundEnd = $i
        break
    }
}
if ($foundStart -ge 1) {
    $rangeStart = $d.Paragraphs.Item($foundStart).Range.Start
    $rangeEnd = $d.Paragraphs.Item($foundEnd).Range.End
    $delRange = $d.Range($rangeStart, $rangeEnd)
    $delRange.Delete()
} else {
    throw "Could not locate paragraphs to delete"
}

# 5) Replace the "הסקירה..." paragraph text with the new "מי שעקב..." text
Replace-ExactText "הסקירה של היום הולכת להיות קצרה וקלילה. אז הקטע ב-LoRA המקורי היה לא לכייל את כל הפרמטרים של מודל שפה (אחרי אימון מקדים) אלא לאמן תוספת לפרמטרים של המודל. כלומר לוקחים את כל הפרמטרים של המודל אחרי המאומן ומאמנים תוספת אליהם שהיא מוגדרת בצורה A*B כאשר A ו- B הם מטריצות בעלות דרגה (ראנק) נמוך (קטנות יותר). " "מי שעקב אחרי המהפכות שהתרחשו בעולם של מודלי שפה ענקיים (LLMs) בטח שמעו על RLHF שזה ראשי תיבות של Reinforcement Learning with Human Feedback (או בקצרה RLHF) בהקשר של אימון מודלי שפה. "

# 6) Replace the "לאחר שמאמנים..." paragraph text with the new "המאמר שנסקור..." text
Replace-ExactText "לאחר שמאמנים מודל שפה (מוצאים ערכים אופטימליים של A ו-B) על דאטהסט ולידציה ומבצעים קווינטוט של המטריצה המקורית ושל התוספת. אז המאמר המסוקר מציע להקטין עוד יותר את מספר הפרמטרים במטריצת התוספת ולהציג אותה כמכפלה של bAdB כאשר מטריצות A ו-B הן קבועות לכל השכבות(נגדמות מהתפלגות נורמלית) וקטורים (לא מטריצות!) b ו-d נלמדות פר שכבה. כך מספר הפרמטרים המנלמדים יורד בצורה משמעות בלי לפגוע בביצועי המודל. בקיצור מודיפיקציה נחמדה של LoRA." "המאמר שנסקור היום ב-#shorthebrewpapereviews לוקח אחת מאבני הבניין של RLHF שזה מודל תגמול (reward model) ומשתמש בה לגנרוט של טקסט. מודל תגמול מיועד לשערוך של איכות הטקסט המגונרט על ידי המודל ומרטת RLHF היא למקסם את התגמול (יחד עם עוד כמה מדדים) במטרה לשפר את איכות הטקסט המגונרט. המאמר המסוקר משתמש למודל התגמול לגנרוט של טקסט בפרט ל`"כיול״ של הסתברויות של הטוקנים שמודל שפה מחשב בשביל לחזות כל טוקן. "

# 7) Insert the three new paragraphs right after the paragraph from step 6
$new_para8 = "המאמר שנסקור היום ב-#shorthebrewpapereviews לוקח אחת מאבני הבניין של RLHF שזה מודל תגמול (reward model) ומשתמש בה לגנרוט של טקסט. מודל תגמול מיועד לשערוך של איכות הטקסט המגונרט על ידי המודל ומרטת RLHF היא למקסם את התגמול (יחד עם עוד כמה מדדים) במטרה לשפר את איכות הטקסט המגונרט. המאמר המסוקר משתמש למודל התגמול לגנרוט של טקסט בפרט ל`"כיול״ של הסתברויות של הטוקנים שמודל שפה מחשב בשביל לחזות כל טוקן. "
$new_para9 = "כלומר עבור כל טוקן נחזה הסתברותו מוזזת בהתאם לתגמול המצופה על ידי הוספת טוקן זה לטוקנים שכבר גונרטו על ידי המודל. טוקנים בעלי הסתברות גבוהה לפי מודל השפה וגם בעלי ערך גבוה של פונקציית התגמול (הממודלת על ידי מודל תגמול) יקבלו עדיפות על פני הטוקנים בעלי ערכי התגמול נמוכים יותר. "
$new_para11 = "מודל התגמול מאומן התאם למשימה נתונה עם פונקציית לוס של המחשבת מרחק בין את התגמול ה-ground truth לזה של המודל לכל טוקן. מעניין כי ככל הקנס על תגמול לא מדויק עולה ככל שהטוקן רחוק יותר מהתחלת הטקסט המגונרט (הקנס על אי דיוק של הטוקן האחרון הוא מקסימלי)."
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text.TrimEnd([char]13, [char]7) -eq $new_para8) {
        $anchorIdx = $i
        break
    }
}
if ($anchorIdx -eq -1) {
    throw "Could not locate anchor paragraph for insertion"
}

$anchor = $d.Paragraphs.Item($anchorIdx)
$anchor.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item($anchorIdx + 1)
$p9.Range.Text = $new_para9

$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Item($anchorIdx + 2)

$p10.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs.Item($anchorIdx + 3)
$p11.Range.Text = $new_para11

Write-Output "done"
